# Update the "想去人数" (interest/attendance count) figures in column F
# for the "展览" (Exhibition) and "全部类型" (All Types) sheets, as produced
# by the site's scheduled data refresh (commit a3196b5).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 941
$ws1.Range("F4").Value  = 115
$ws1.Range("F6").Value  = 40
$ws1.Range("F7").Value  = 6401
$ws1.Range("F9").Value  = 893
$ws1.Range("F10").Value = 249
$ws1.Range("F11").Value = 729
$ws1.Range("F12").Value = 478
$ws1.Range("F14").Value = 29
$ws1.Range("F15").Value = 360
$ws1.Range("F16").Value = 809
$ws1.Range("F17").Value = 2290
$ws1.Range("F18").Value = 39
$ws1.Range("F19").Value = 148
$ws1.Range("F20").Value = 726
$ws1.Range("F21").Value = 28
$ws1.Range("F22").Value = 384
$ws1.Range("F23").Value = 170
$ws1.Range("F24").Value = 48
$ws1.Range("F25").Value = 64
$ws1.Range("F26").Value = 3
$ws1.Range("F27").Value = 214

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 941
$ws4.Range("F5").Value  = 115
$ws4.Range("F10").Value = 40
$ws4.Range("F11").Value = 6401
$ws4.Range("F13").Value = 893
$ws4.Range("F14").Value = 249
$ws4.Range("F15").Value = 729
$ws4.Range("F16").Value = 478
$ws4.Range("F18").Value = 29
$ws4.Range("F19").Value = 360
$ws4.Range("F20").Value = 809
$ws4.Range("F22").Value = 2290
$ws4.Range("F23").Value = 39
$ws4.Range("F25").Value = 148
$ws4.Range("F26").Value = 726
$ws4.Range("F27").Value = 28
$ws4.Range("F28").Value = 384
$ws4.Range("F29").Value = 170
$ws4.Range("F30").Value = 48
$ws4.Range("F31").Value = 64
$ws4.Range("F32").Value = 3
$ws4.Range("F33").Value = 214
